# Applies the "Updated cryptos list" data refresh described by the commit diff.
# For each changed cell we assign the new literal text. Numeric-looking strings
# (e.g. "680.86") are prefixed with a leading apostrophe so Excel stores them as
# text (matching the original inlineStr cell type) instead of coercing to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '69.478.13'
$ws.Range("E2").Value = '  +0.07%  '
# Row 3
$ws.Range("D3").Value = '3.690.83'
$ws.Range("E3").Value = '  -0.01%  '
# Row 4
$ws.Range("E4").Value = '  +0.02%  '
# Row 5
$ws.Range("D5").Value = '''680.86'
$ws.Range("E5").Value = '  -0.83%  '
# Row 6
$ws.Range("D6").Value = '''161.46'
$ws.Range("E6").Value = '  +0.78%  '
# Row 7
$ws.Range("E7").Value = '  -0.08%  '
# Row 8
$ws.Range("D8").Value = '''0.495'
$ws.Range("E8").Value = '  +0.19%  '
# Row 9
$ws.Range("E9").Value = '  +0.20%  '
# Row 10
$ws.Range("D10").Value = '''7.14'
$ws.Range("E10").Value = '  -0.96%  '
# Row 11
$ws.Range("E11").Value = '  +0.64%  '
# Row 12
$ws.Range("D12").Value = '''0.0000233'
$ws.Range("E12").Value = '  +0.10%  '
# Row 13
$ws.Range("D13").Value = '4.316.17'
$ws.Range("E13").Value = '  +0.04%  '
# Row 14
$ws.Range("D14").Value = '''32.46'
$ws.Range("E14").Value = '  -0.18%  '
# Row 15
$ws.Range("D15").Value = '3.700.03'
$ws.Range("E15").Value = '  +0.25%  '
# Row 16
$ws.Range("D16").Value = '69.486.15'
$ws.Range("E16").Value = '  +0.09%  '
# Row 17
$ws.Range("E17").Value = '  +2.42%  '
# Row 18
$ws.Range("D18").Value = '''16.02'
$ws.Range("E18").Value = '  +0.39%  '
# Row 19
$ws.Range("D19").Value = '''6.48'
$ws.Range("E19").Value = '  +0.47%  '
# Row 20
$ws.Range("D20").Value = '''471.10'
$ws.Range("E20").Value = '  -0.09%  '
# Row 21
$ws.Range("E21").Value = '  -1.57%  '
# Row 22
$ws.Range("D22").Value = '''0.650'
# Row 23
$ws.Range("D23").Value = '''80.48'
# Row 24
$ws.Range("D24").Value = '3.838.87'
$ws.Range("E24").Value = '  -0.03%  '
# Row 25
$ws.Range("E25").Value = '  -0.05%  '
# Row 26
$ws.Range("D26").Value = '''0.0000126'
$ws.Range("E26").Value = '  +0.77%  '
# Row 27
$ws.Range("D27").Value = '''10.89'
$ws.Range("E27").Value = '  -1.20%  '
# Row 28
$ws.Range("D28").Value = '''9.12'
$ws.Range("E28").Value = '  -1.02%  '
# Row 29
$ws.Range("E29").Value = '  +0.03%  '
# Row 30
$ws.Range("D30").Value = '''1.74'
$ws.Range("E30").Value = '  -0.72%  '
# Row 31
$ws.Range("E31").Value = '  -0.27%  '
# Row 32
$ws.Range("D32").Value = '''6.57'
$ws.Range("E32").Value = '  -1.41%  '
# Row 33
$ws.Range("E33").Value = '  +0.44%  '
# Row 34
$ws.Range("D34").Value = '''27.00'
$ws.Range("E34").Value = '  +0.94%  '
# Row 35
$ws.Range("D35").Value = '3.682.20'
# Row 36
$ws.Range("E36").Value = '  +1.66%  '
# Row 37
$ws.Range("D37").Value = '''8.46'
$ws.Range("E37").Value = '  +3.20%  '
# Row 38
$ws.Range("D38").Value = '''6.21'
$ws.Range("E38").Value = '  +0.82%  '
# Row 39
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").Value = '''2.28'
$ws.Range("E39").Value = '  -0.12%  '
# Row 40
$ws.Range("B40").Value = 'USDe'
$ws.Range("C40").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D40").Value = '''1.00'
$ws.Range("E40").Value = '  +0.00%  '
# Row 41
$ws.Range("E41").Value = '  +0.02%  '
# Row 42
$ws.Range("D42").Value = '''0.0900'
$ws.Range("E42").Value = '  -0.50%  '
# Row 43
$ws.Range("D43").Value = '''168.63'
$ws.Range("E43").Value = '  +1.18%  '
# Row 44
$ws.Range("D44").Value = '''0.942'
$ws.Range("E44").Value = '  -0.11%  '
# Row 45
$ws.Range("D45").Value = '''46.76'
$ws.Range("E45").Value = '  -2.34%  '
# Row 46
$ws.Range("D46").Value = '''2.74'
$ws.Range("E46").Value = '  +0.38%  '
# Row 47
$ws.Range("D47").Value = '''0.000279'
$ws.Range("E47").Value = '  +1.37%  '
# Row 48
$ws.Range("D48").Value = '''1.29'
$ws.Range("E48").Value = '  -0.97%  '
# Row 49
$ws.Range("D49").Value = '''27.60'
$ws.Range("E49").Value = '  -2.93%  '
# Row 50
$ws.Range("E50").Value = '  -2.71%  '
# Row 51
$ws.Range("D51").Value = '''7.90'
$ws.Range("E51").Value = '  +1.08%  '

Write-Output "Applied 88 cell updates across 49 rows"
